# Update "想去人数" (attendance count) figures in the 展览 (F column) sheet
# and the matching rows in the 全部类型 (All types) aggregate sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 455
$wsExpo.Range("F3").Value = 6
$wsExpo.Range("F5").Value = 78
$wsExpo.Range("F6").Value = 5242
$wsExpo.Range("F8").Value = 86
$wsExpo.Range("F9").Value = 99
$wsExpo.Range("F10").Value = 353
$wsExpo.Range("F11").Value = 12
$wsExpo.Range("F12").Value = 65

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 455
$wsAll.Range("F4").Value = 6
$wsAll.Range("F9").Value = 78
$wsAll.Range("F10").Value = 5242
$wsAll.Range("F12").Value = 86
$wsAll.Range("F13").Value = 99
$wsAll.Range("F15").Value = 353
$wsAll.Range("F16").Value = 12
$wsAll.Range("F17").Value = 65
